# "A bit of cleaning and analytics"
#
# Re-labels the Sheet1 header row, drops the two now-unused trailing
# columns (E/F), and backfills a "file format" column (PDF/CSV/TXT/JSON)
# for the years that have electronic reports, with the new data cells
# centered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename/repurpose the report columns, drop E:F ---
$ws.Range("B1").Value = "Score Table"
$ws.Range("C1").Value = "Country Rankings"
$ws.Range("D1").Value = "Awards"
$ws.Range("E1:F1").Clear()

# --- Backfill the "Awards"/export-format data, year by year ---
# 2005-2009: PDF only
$ws.Range("D18:D22").Value = "PDF"

# 2010-2013: CSV + PDF
$ws.Range("C23:C26").Value = "CSV"
$ws.Range("D23:D26").Value = "PDF"

# 2014-2015: CSV + TXT
$ws.Range("C27:C28").Value = "CSV"
$ws.Range("D27:D28").Value = "TXT"

# 2016-2019: CSV + CSV + JSON
$ws.Range("B29:B32").Value = "CSV"
$ws.Range("C29:C32").Value = "CSV"
$ws.Range("D29:D32").Value = "JSON"

# --- Center-align the newly populated data block ---
$ws.Range("B2:D32").HorizontalAlignment = -4108

# --- Column widths: B/C/D resized, A reverts to default ---
$ws.Columns.Item(2).ColumnWidth = 14.0
$ws.Columns.Item(3).ColumnWidth = 19.0
$ws.Columns.Item(4).ColumnWidth = 9.833333333333334

# --- Selection moves to the last newly-touched cell ---
$ws.Range("D18").Select() | Out-Null
